{"js": "// Replace the date and each two-digit multiplication problem in the\n// document with its updated value, preserving all run/paragraph formatting.\n// Each old value is unique in the document, so a single exact-match\n// search + in-place \"Replace\" insertText per pair is safe and order\n// independent.\nconst replacements = [\n  [\"2025-11-12 Wednesday\", \"2025-11-13 Thursday\"],\n  [\"50\u00d721=\", \"65\u00d787=\"],\n  [\"24\u00d770=\", \"12\u00d777=\"],\n  [\"93\u00d731=\", \"60\u00d758=\"],\n  [\"89\u00d780=\", \"59\u00d796=\"],\n  [\"68\u00d767=\", \"80\u00d718=\"],\n  [\"58\u00d758=\", \"42\u00d729=\"],\n  [\"45\u00d780=\", \"58\u00d715=\"],\n  [\"23\u00d721=\", \"22\u00d749=\"],\n  [\"55\u00d736=\", \"78\u00d758=\"],\n  [\"59\u00d740=\", \"77\u00d720=\"],\n  [\"71\u00d715=\", \"93\u00d781=\"],\n  [\"66\u00d739=\", \"93\u00d716=\"],\n  [\"98\u00d764=\", \"56\u00d795=\"],\n  [\"31\u00d781=\", \"98\u00d762=\"],\n  [\"18\u00d726=\", \"90\u00d788=\"],\n  [\"97\u00d768=\", \"91\u00d717=\"],\n  [\"95\u00d763=\", \"92\u00d776=\"],\n  [\"63\u00d732=\", \"84\u00d755=\"],\n  [\"73\u00d717=\", \"91\u00d725=\"],\n  [\"86\u00d728=\", \"64\u00d724=\"],\n  [\"45\u00d739=\", \"15\u00d727=\"],\n  [\"97\u00d779=\", \"42\u00d720=\"],\n  [\"31\u00d751=\", \"53\u00d753=\"],\n  [\"90\u00d764=\", \"51\u00d717=\"],\n  [\"34\u00d773=\", \"39\u00d741=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date and each two-digit multiplication problem in the\n# document with its updated value, preserving all run/paragraph formatting.\n# Each old value is unique in the document, so a single exact-match\n# Find/Replace pass per pair is safe and order independent.\n\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2025-11-12 Wednesday\", \"2025-11-13 Thursday\"),\n    @(\"50\u00d721=\", \"65\u00d787=\"),\n    @(\"24\u00d770=\", \"12\u00d777=\"),\n    @(\"93\u00d731=\", \"60\u00d758=\"),\n    @(\"89\u00d780=\", \"59\u00d796=\"),\n    @(\"68\u00d767=\", \"80\u00d718=\"),\n    @(\"58\u00d758=\", \"42\u00d729=\"),\n    @(\"45\u00d780=\", \"58\u00d715=\"),\n    @(\"23\u00d721=\", \"22\u00d749=\"),\n    @(\"55\u00d736=\", \"78\u00d758=\"),\n    @(\"59\u00d740=\", \"77\u00d720=\"),\n    @(\"71\u00d715=\", \"93\u00d781=\"),\n    @(\"66\u00d739=\", \"93\u00d716=\"),\n    @(\"98\u00d764=\", \"56\u00d795=\"),\n    @(\"31\u00d781=\", \"98\u00d762=\"),\n    @(\"18\u00d726=\", \"90\u00d788=\"),\n    @(\"97\u00d768=\", \"91\u00d717=\"),\n    @(\"95\u00d763=\", \"92\u00d776=\"),\n    @(\"63\u00d732=\", \"84\u00d755=\"),\n    @(\"73\u00d717=\", \"91\u00d725=\"),\n    @(\"86\u00d728=\", \"64\u00d724=\"),\n    @(\"45\u00d739=\", \"15\u00d727=\"),\n    @(\"97\u00d779=\", \"42\u00d720=\"),\n    @(\"31\u00d751=\", \"53\u00d753=\"),\n    @(\"90\u00d764=\", \"51\u00d717=\"),\n    @(\"34\u00d773=\", \"39\u00d741=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $rng = $d.Content\n    $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
